$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Targeted scalar fixes -------------------------------------------------
# "isPivot" flag recomputed for the week of 2024-06-03 (row 1146).
$ws.Range("O1146").Value = 2

# "detect_structure" recomputed to 0 for these early (2003) weeks.
$ws.Range("Q55").Value = 0
$ws.Range("Q61").Value = 0
$ws.Range("Q65").Value = 0

# "backup" column backfilled with 0 for the last two previously-blank rows.
$ws.Range("R1148").Value = 0
$ws.Range("R1149").Value = 0

# --- Append newly completed weekly bars (2024-07-01 .. 2024-12-23) --------
# CSV layout: Row,A,B,C,D,E,F,G,H,I,J,K,L,M,N,O,P,Q,R
# Columns F (Adj Close) and R (backup) are intentionally left blank for
# these brand-new rows, matching the source export (not yet back-filled).
$newRows = @"
1150,45474,7001.39990234375,7127.75,6980,7054.4501953125,,405405,2024,7,1,0,0,0,27,0,0,0,
1151,45481,7061,7175,6931.35009765625,7065.5498046875,,1073776,2024,7,8,0,0,0,28,0,0,0,
1152,45488,7079,7103.7998046875,6942,6980.5498046875,,362126,2024,7,15,0,0,0,29,0,0,1,
1153,45495,6978.9501953125,7069.14990234375,6852,6961.60009765625,,453784,2024,7,22,0,0,0,30,0,0,0,
1154,45502,6975,7011.85009765625,6890.0498046875,6911.0498046875,,304300,2024,7,29,0,0,0,31,0,0,0,
1155,45509,6831.2001953125,6899.89990234375,6666,6875.75,,358940,2024,8,5,0,0,0,32,2,0,0,
1156,45516,6875.75,6888,6724.64990234375,6799.64990234375,,234065,2024,8,12,0,0,0,33,0,0,0,
1157,45523,6799.64990234375,7139,6780,7088.14990234375,,546618,2024,8,19,0,0,0,34,0,0,0,
1158,45530,7149,9080,7103,8000.7001953125,,12098917,2024,8,26,0,0,0,35,1,0,0,
1159,45537,8009,8050,7567.10009765625,7589.4501953125,,1248625,2024,9,2,0,0,0,36,0,0,0,
1160,45544,7590,8024,7490,7714.89990234375,,1490515,2024,9,9,0,0,0,37,0,0,0,
1161,45551,7750,7939,7496.7001953125,7859.39990234375,,818530,2024,9,16,0,0,0,38,0,0,1,
1162,45558,7885.0498046875,8142.14990234375,7600,7811.39990234375,,1107017,2024,9,23,0,0,0,39,0,0,0,
1163,45565,7780,7807.9501953125,7451,7576.10009765625,,476841,2024,9,30,0,0,0,40,0,0,0,
1164,45572,7580,8027,7352,7675.14990234375,,1624637,2024,10,7,0,0,0,41,0,0,2,
1165,45579,7728.39990234375,7749,7275,7328.7001953125,,674749,2024,10,14,0,0,0,42,0,0,0,
1166,45586,7365,7485,6870.60009765625,6990.64990234375,,773240,2024,10,21,0,0,0,43,0,0,0,
1167,45593,7068,7134,6969,7107.5,,317745,2024,10,28,0,0,0,44,0,0,0,
1168,45600,7146.9501953125,7390,6821.64990234375,6940.7998046875,,464521,2024,11,4,0,0,0,45,0,0,0,
1169,45607,6940.7998046875,7011.0498046875,6338.2001953125,6374.5498046875,,1454666,2024,11,11,0,0,0,46,0,0,0,
1170,45614,6405.0498046875,6620,6286,6586.4501953125,,589771,2024,11,18,0,0,0,47,2,0,0,
1171,45621,6690,6879,6610.14990234375,6670.75,,706862,2024,11,25,0,0,0,48,0,0,0,
1172,45628,6670.75,7474,6613,7396.4501953125,,2148940,2024,12,2,0,0,0,49,0,0,0,
1173,45635,7400,7460,7211.9501953125,7364.39990234375,,819970,2024,12,9,0,0,0,50,0,0,0,
1174,45642,7372.5,7424.35009765625,6900,6921.60009765625,,765032,2024,12,16,0,0,0,51,0,0,0,
1175,45649,7037.7001953125,7037.7001953125,6764.39990234375,6886.35009765625,,476300,2024,12,23,0,0,0,52,0,0,0,
"@

$lines = $newRows -split "`n"
foreach ($line in $lines) {
    $line = $line.Trim()
    if ($line -eq "") { continue }

    $parts = $line.Split(",")
    $r = [int]$parts[0]

    for ($c = 1; $c -le 18; $c++) {
        $cellText = $parts[$c]
        if ($cellText -ne "") {
            $ws.Cells.Item($r, $c).Value = [double]$cellText
        }
    }
}

# Column A holds Datetime serials; reuse the workbook's existing custom
# date/time format so the new cells pick up the same style as the rest of
# column A (s="2" in the OOXML).
$ws.Range("A1150:A1175").NumberFormat = "YYYY-MM-DD HH:MM:SS"
